# This script expands the abbreviated party codes used throughout the
# workbook's shared strings into their full "Code - Name (Name, Code)"
# form, e.g. "BQ " -> "BQ - Bloc Québécois (Bloc Québécois, BQ)".
#
# The party codes appear both on their own (e.g. in the "Transformed Data"
# / "Parties per Year" sheets) and combined into "+"-joined coalition
# labels (e.g. "BQ +Ind") across several other sheets (Coalitions,
# Winning/Minimal Winning/Maximal Losing/Unique Tying Coalitions, etc.),
# as well as inside Python-tuple-style strings such as
# "('Lib+ND ', 'BQ +Con+Non')".
#
# Using Find/Replace with LookAt:=xlPart (partial match) across every
# worksheet updates every occurrence - standalone codes and all the
# "+"-joined combinations alike - in a single pass per code, without
# needing to enumerate every individual cell or combination.
#
# IMPORTANT: this runtime's Range.Replace ignores the MatchCase argument
# and always matches case-insensitively. That is a problem because some
# of the short codes are case-insensitive substrings of other codes'
# *expanded* replacement text (e.g. "nd " occurs inside "Ind - ..." and
# would collide with the "ND " code; "Con" occurs inside "...Progressive
# Conservative..." and would collide with the "Con" code). Doing the 9
# replacements directly, one after another, therefore corrupts text that
# was already expanded by an earlier pass.
#
# To avoid this, expansion is done in two phases:
#   1) Replace each short code with a unique, collision-free placeholder
#      token. Because all 9 placeholders are distinct and do not overlap
#      with each other, with any of the short codes, or with any existing
#      text in the workbook, this phase is safe even under
#      case-insensitive matching.
#   2) Replace each placeholder token with its final expanded text. Since
#      the placeholders are unique literal tokens, this phase cannot
#      accidentally re-match text inserted by earlier steps either.

$wb = $excel.ActiveWorkbook

$codes = @(
    @("BQ ", "BQ - Bloc Québécois (Bloc Québécois, BQ)"),
    @("Ind", "Ind - Independents (Independents, Ind)"),
    @("Lib", "Lib - Liberal  (Liberal , Lib)"),
    @("ND ", "ND - New Democratic (New Democratic, ND)"),
    @("PC ", "PC - Progressive Conservative  (Progressive Conservative , PC)"),
    @("RPC", "RPC/RP - Reform Party of Canada / Canadian Alliance (Reform Party of Canada / Canadian Alliance, RPC/RP)"),
    @("Non", "None - No-Affiliation (No-Affiliation, None)"),
    @("Con", "Con - Conservative (Conservative, Con)"),
    @("GP ", "GP - Green Party of Canada (Green Party of Canada, GP)")
)

# Phase 1: short code -> unique placeholder token
for ($i = 0; $i -lt $codes.Count; $i++) {
    $what = $codes[$i][0]
    $placeholder = "@@@PARTYPLACEHOLDER" + $i + "@@@"
    foreach ($ws in $wb.Worksheets) {
        $ws.Cells.Replace($what, $placeholder, 2, 1, $false, $true) | Out-Null
    }
}

# Phase 2: placeholder token -> final expanded text
for ($i = 0; $i -lt $codes.Count; $i++) {
    $withText = $codes[$i][1]
    $placeholder = "@@@PARTYPLACEHOLDER" + $i + "@@@"
    foreach ($ws in $wb.Worksheets) {
        $ws.Cells.Replace($placeholder, $withText, 2, 1, $false, $true) | Out-Null
    }
}
